$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateSTP")

# Update the "04011" suffixed test-data values on row 2/3 (Data5/Data6
# columns) to the new "04012" suffix for the List Page category search case.
$ws.Range("F2").Value = "FullName104012"
$ws.Range("F3").Value = "Short104012"
$ws.Range("G2").Value = "FullName204012"
$ws.Range("G3").Value = "Short204012"

# Move the active selection to D16, as in the saved workbook.
$ws.Activate()
$ws.Range("D16").Select()
